# Weekly driver report update for 2025-04-19
# Updates the "Driver Summary" sheet: refreshes the Bad/Good driver tables
# with this week's roaming data (Wi-Fi 6 AX200 family replacing AX201).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Bad Drivers table (rows 3-6)
#    The Realtek row (row 5) drops out of the data this week, so remove
#    it and let everything below shift up - this also slides the "Good
#    Drivers" block up by one row, matching the new layout.
# ---------------------------------------------------------------------
$ws.Rows("5").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# Row 3: Intel AX200 160MHz - 23.110.0.5
$ws.Range("A3").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 23.110.0.5"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 26
$ws.Range("D3").Value = 93.3

# Row 4: Realtek adapter, refreshed counts
$ws.Range("A4").Value = "Realtek RTL8852AE WiFi 6 802.11ax PCIe Adapter - 6001.10.356.0"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 161
$ws.Range("D4").Value = 98.59999999999999

# Row 5: Totals
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 187

# ---------------------------------------------------------------------
# 2) Good Drivers table (now starting at row 11 after the shift above)
#    Replace the six existing AX201 rows with the new AX200 rows, then
#    stamp seven more rows (copying the formatting of the last templated
#    row down) and fill them with the rest of this week's good drivers.
# ---------------------------------------------------------------------

# Grow the table: rows 19-25 need the same formatting as row 18 already
# has (bold/right-aligned number styles etc.) before we fill in values.
$ws.Range("A18:E18").Copy($ws.Range("A19"))
$ws.Range("A18:E18").Copy($ws.Range("A20"))
$ws.Range("A18:E18").Copy($ws.Range("A21"))
$ws.Range("A18:E18").Copy($ws.Range("A22"))
$ws.Range("A18:E18").Copy($ws.Range("A23"))
$ws.Range("A18:E18").Copy($ws.Range("A24"))
$ws.Range("A18:E18").Copy($ws.Range("A25"))

# Row 13
$ws.Range("A13").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.160.0.4"
$ws.Range("B13").Value = 96526
$ws.Range("C13").ClearContents()
$ws.Range("D13").Value = 99.90000000000001
$ws.Range("E13").ClearContents()

# Row 14
$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.230.0.8"
$ws.Range("B14").Value = 328411
$ws.Range("C14").ClearContents()
$ws.Range("D14").Value = 99.90000000000001
$ws.Range("E14").ClearContents()

# Row 15
$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.200.0.6"
$ws.Range("B15").Value = 143808
$ws.Range("C15").ClearContents()
$ws.Range("D15").Value = 99.90000000000001
$ws.Range("E15").ClearContents()

# Row 16
$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.190.0.4"
$ws.Range("B16").Value = 287148
$ws.Range("C16").ClearContents()
$ws.Range("D16").Value = 99.90000000000001
$ws.Range("E16").ClearContents()

# Row 17
$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.250.10.1"
$ws.Range("B17").Value = 69578
$ws.Range("C17").ClearContents()
$ws.Range("D17").Value = 99.90000000000001
$ws.Range("E17").ClearContents()

# Row 18
$ws.Range("A18").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.30.0.11"
$ws.Range("B18").Value = 67111
$ws.Range("C18").ClearContents()
$ws.Range("D18").Value = 100
$ws.Range("E18").ClearContents()

# Row 19 (new)
$ws.Range("A19").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 21.30.4.1"
$ws.Range("B19").Value = 13016
$ws.Range("C19").ClearContents()
$ws.Range("D19").Value = 100
$ws.Range("E19").ClearContents()

# Row 20 (new)
$ws.Range("A20").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 23.70.2.3"
$ws.Range("B20").Value = 18721
$ws.Range("C20").ClearContents()
$ws.Range("D20").Value = 99.90000000000001
$ws.Range("E20").Value = "2024-07-23"

# Row 21 (new)
$ws.Range("A21").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.10.0.7"
$ws.Range("B21").Value = 66577
$ws.Range("C21").ClearContents()
$ws.Range("D21").Value = 100
$ws.Range("E21").Value = "2024-05-09"

# Row 22 (new)
$ws.Range("A22").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 21.60.2.1"
$ws.Range("B22").Value = 26241
$ws.Range("C22").ClearContents()
$ws.Range("D22").Value = 100
$ws.Range("E22").Value = "2021-01-19"

# Row 23 (new)
$ws.Range("A23").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.0.1.1"
$ws.Range("B23").Value = 15730
$ws.Range("C23").ClearContents()
$ws.Range("D23").Value = 99.90000000000001
$ws.Range("E23").Value = "2020-09-28"

# Row 24 (new)
$ws.Range("A24").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 21.40.2.2"
$ws.Range("B24").Value = 88435
$ws.Range("C24").ClearContents()
$ws.Range("D24").Value = 99.90000000000001
$ws.Range("E24").Value = "2019-08-31"

# Row 25 (new)
$ws.Range("A25").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 21.10.1.2"
$ws.Range("B25").Value = 46270
$ws.Range("C25").ClearContents()
$ws.Range("D25").Value = 100
$ws.Range("E25").Value = "2019-04-23"
